$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a single new row before row 15. The two previously-blank rows
# (13 and 14) shrink to one blank row after the two new content rows are
# added, so the net shift for everything at/after row 15 is +1.
$ws.Rows("15:15").Insert()

# New row 13: Komiku track credit
$ws.Range("A13").Value = "Komiku_-_02_-_Boss_4__Cobblestone_in_their_face.mp3"
$ws.Range("B13").Value = "https://www.chosic.com/download-audio/25453/"
$ws.Range("C13").Value = "You are free to use this music in your projects with no required crediting. However, linking back is greatly appreciated. You can use the following text"
$ws.Range("C13").Style = $ws.Range("C11").Style

# New row 14: Loyalty Freak Music track credit
$ws.Range("A14").Value = "Loyalty_Freak_Music_-_04_-_Cant_Stop_My_Feet_.mp3"
$ws.Range("B14").Value = "https://www.chosic.com/download-audio/25495/"
$ws.Range("C14").Value = "You are free to use this music in your projects with no required crediting. However, linking back is greatly appreciated. You can use the following text"
$ws.Range("C14").Style = $ws.Range("C11").Style

# The hyperlink collection does not automatically re-anchor itself to the
# shifted rows (B16/B17 stay put instead of moving to B17/B18), so rebuild
# the hyperlinks collection from scratch in the same rId order.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B6"), "https://opengameart.org/content/spikes-0")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.fontspace.com/a-area-kilometer-50-font-f53888")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://freesound.org/people/Whiprealgood/sounds/87535/")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://freesound.org/people/suntemple/sounds/253172/")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://opengameart.org/content/simple-explosion-bleeds-game-art")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://opengameart.org/content/various-inventory-24-pixel-icon-set")
$ws.Hyperlinks.Add($ws.Range("B17"), "https://elthen.itch.io/2d-pixel-art-vegetable-monsters-sprite-pack")
$ws.Hyperlinks.Add($ws.Range("B18"), "https://free-game-assets.itch.io/night-city-street-2d-background-tiles")
$ws.Hyperlinks.Add($ws.Range("B8"), "https://opengameart.org/content/energy-icon")
